$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# Values are assigned in an order that reproduces the target shared-string
# table order: email, nom, civilite, tel, adresse, prenom, formation
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "nom"
$ws.Range("D1").Value = "civilite"
$ws.Range("E1").Value = "tel"
$ws.Range("F1").Value = "adresse"
$ws.Range("C1").Value = "prenom"
$ws.Range("G1").Value = "formation"

# --- Add row 2 with A2 formatted using the built-in Hyperlink style ---
# A temporary hyperlink is created (this registers the built-in
# "Hyperlink" cell style/font in the workbook), then the link and its
# text are removed again, leaving an empty cell that still carries the
# hyperlink style.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:test@test.com") | Out-Null
$ws.Hyperlinks.Delete() | Out-Null
$ws.Range("A2").ClearContents() | Out-Null

# --- Update the selected cell shown when the workbook is reopened ---
$ws.Range("E20").Select() | Out-Null
